$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = 151
$ws.Cells.Item(2, 5).Value = 16.5
$ws.Cells.Item(3, 4).Value = 107
$ws.Cells.Item(4, 4).Value = 89
$ws.Cells.Item(5, 4).Value = 172
$ws.Cells.Item(7, 4).Value = 40
$ws.Cells.Item(9, 4).Value = 147
$ws.Cells.Item(10, 4).Value = 149
$ws.Cells.Item(11, 4).Value = 95
$ws.Cells.Item(12, 4).Value = 78
$ws.Cells.Item(13, 4).Value = 182
$ws.Cells.Item(14, 4).Value = 43
$ws.Cells.Item(15, 4).Value = 32
$ws.Cells.Item(16, 4).Value = 70
$ws.Cells.Item(17, 4).Value = 154
$ws.Cells.Item(18, 4).Value = 155
$ws.Cells.Item(19, 4).Value = 97
$ws.Cells.Item(20, 4).Value = 73
$ws.Cells.Item(21, 4).Value = 187
$ws.Cells.Item(22, 4).Value = 44
$ws.Cells.Item(23, 4).Value = 40
$ws.Cells.Item(24, 4).Value = 81
$ws.Cells.Item(25, 4).Value = 160
$ws.Cells.Item(26, 4).Value = 159
$ws.Cells.Item(27, 4).Value = 86
$ws.Cells.Item(28, 4).Value = 79
$ws.Cells.Item(29, 4).Value = 184
$ws.Cells.Item(30, 4).Value = 40
$ws.Cells.Item(32, 4).Value = 61
$ws.Cells.Item(33, 4).Value = 153
$ws.Cells.Item(34, 4).Value = 143
$ws.Cells.Item(35, 4).Value = 73
$ws.Cells.Item(36, 4).Value = 68
$ws.Cells.Item(37, 4).Value = 181
$ws.Cells.Item(38, 4).Value = 38
$ws.Cells.Item(39, 4).Value = 45
$ws.Cells.Item(40, 4).Value = 64
$ws.Cells.Item(41, 4).Value = 150
$ws.Cells.Item(41, 5).Value = 17.7
$ws.Cells.Item(42, 4).Value = 134
$ws.Cells.Item(43, 4).Value = 78
$ws.Cells.Item(44, 4).Value = 68
$ws.Cells.Item(45, 4).Value = 187
$ws.Cells.Item(46, 4).Value = 38
$ws.Cells.Item(47, 4).Value = 44
$ws.Cells.Item(48, 4).Value = 75
$ws.Cells.Item(49, 4).Value = 148
$ws.Cells.Item(50, 4).Value = 116
$ws.Cells.Item(51, 4).Value = 86
$ws.Cells.Item(52, 4).Value = 63
$ws.Cells.Item(53, 4).Value = 184
$ws.Cells.Item(54, 4).Value = 47
$ws.Cells.Item(54, 5).Value = 5.3
$ws.Cells.Item(55, 4).Value = 44
$ws.Cells.Item(56, 4).Value = 79
$ws.Cells.Item(57, 4).Value = 152
$ws.Cells.Item(58, 4).Value = 124
$ws.Cells.Item(59, 4).Value = 80
$ws.Cells.Item(60, 4).Value = 74
$ws.Cells.Item(61, 4).Value = 185
$ws.Cells.Item(62, 4).Value = 45
$ws.Cells.Item(63, 4).Value = 51
$ws.Cells.Item(64, 4).Value = 89
$ws.Cells.Item(64, 5).Value = 9.699999999999999
$ws.Cells.Item(65, 4).Value = 165
$ws.Cells.Item(65, 5).Value = 17.8
$ws.Cells.Item(66, 4).Value = 137
$ws.Cells.Item(66, 5).Value = 16.4
$ws.Cells.Item(67, 4).Value = 71
$ws.Cells.Item(67, 5).Value = 8.6
$ws.Cells.Item(68, 4).Value = 56
$ws.Cells.Item(68, 5).Value = 6.8
$ws.Cells.Item(69, 4).Value = 154
$ws.Cells.Item(70, 4).Value = 37
$ws.Cells.Item(71, 4).Value = 45
$ws.Cells.Item(72, 4).Value = 95
$ws.Cells.Item(72, 5).Value = 11.4
$ws.Cells.Item(73, 4).Value = 165
$ws.Cells.Item(73, 5).Value = 19.8
$ws.Cells.Item(74, 4).Value = 143
$ws.Cells.Item(74, 5).Value = 15.2
$ws.Cells.Item(75, 4).Value = 82
$ws.Cells.Item(76, 4).Value = 65
$ws.Cells.Item(76, 5).Value = 6.9
$ws.Cells.Item(77, 4).Value = 202
$ws.Cells.Item(78, 4).Value = 40
$ws.Cells.Item(79, 4).Value = 41
$ws.Cells.Item(80, 4).Value = 88
$ws.Cells.Item(80, 5).Value = 9.4
$ws.Cells.Item(81, 4).Value = 169
$ws.Cells.Item(81, 5).Value = 18
$ws.Cells.Item(82, 4).Value = 121
$ws.Cells.Item(82, 5).Value = 13.1
$ws.Cells.Item(83, 4).Value = 76
$ws.Cells.Item(84, 4).Value = 74
$ws.Cells.Item(84, 5).Value = 8
$ws.Cells.Item(85, 4).Value = 195
$ws.Cells.Item(86, 4).Value = 38
$ws.Cells.Item(87, 4).Value = 47
$ws.Cells.Item(88, 4).Value = 87
$ws.Cells.Item(89, 4).Value = 182
$ws.Cells.Item(89, 5).Value = 19.6
$ws.Cells.Item(90, 4).Value = 116
$ws.Cells.Item(90, 5).Value = 12.7
$ws.Cells.Item(91, 4).Value = 70
$ws.Cells.Item(92, 4).Value = 71
$ws.Cells.Item(93, 4).Value = 181
$ws.Cells.Item(93, 5).Value = 19.9
$ws.Cells.Item(94, 4).Value = 40
$ws.Cells.Item(94, 5).Value = 4.3
$ws.Cells.Item(95, 4).Value = 48
$ws.Cells.Item(96, 4).Value = 86
$ws.Cells.Item(96, 5).Value = 9.4
$ws.Cells.Item(97, 4).Value = 193
$ws.Cells.Item(97, 5).Value = 21.2
$ws.Cells.Item(98, 4).Value = 111
$ws.Cells.Item(99, 4).Value = 78
$ws.Cells.Item(100, 4).Value = 80
$ws.Cells.Item(100, 5).Value = 8.199999999999999
$ws.Cells.Item(101, 4).Value = 193
$ws.Cells.Item(102, 4).Value = 48
$ws.Cells.Item(103, 4).Value = 56
$ws.Cells.Item(104, 4).Value = 97
$ws.Cells.Item(104, 5).Value = 9.9
$ws.Cells.Item(105, 4).Value = 201
